# Insert a new data row into the sheet at row 434, shifting existing rows
# 434..530 down to 435..531, and populate the new row with the values
# described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 434 (this shifts rows 434:530 down to 435:531)
$ws.Rows.Item(434).Insert()

# Populate the newly inserted row 434 with the new data record
$ws.Cells.Item(434, 1).Value = 5
$ws.Cells.Item(434, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(434, 3).Value = "Maule"
$ws.Cells.Item(434, 4).Value = 45244
$ws.Cells.Item(434, 4).NumberFormat = $ws.Cells.Item(435, 4).NumberFormat
$ws.Cells.Item(434, 5).Value = 7
$ws.Cells.Item(434, 6).Value = 100112008
$ws.Cells.Item(434, 7).Value = "Coliflor"
$ws.Cells.Item(434, 8).Value = "Sin especificar"
$ws.Cells.Item(434, 9).Value = "Primera"
$ws.Cells.Item(434, 10).Value = 3000
$ws.Cells.Item(434, 11).Value = 1000
$ws.Cells.Item(434, 12).Value = 1000
$ws.Cells.Item(434, 13).Value = 1000
$ws.Cells.Item(434, 14).Value = "`$/unidad"
$ws.Cells.Item(434, 15).Value = "Región del Maule"
$ws.Cells.Item(434, 16).Value = 1000
$ws.Cells.Item(434, 17).Value = 1
$ws.Cells.Item(434, 18).Value = "Hortaliza"
